$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D can look like plain numbers (e.g. "91.89"), which Excel
# would otherwise auto-convert to a numeric value on assignment. Force them to
# be stored as text (matching the original inlineStr cells), then restore the
# cell style so no permanent number-format change is left behind.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "41.340.97"
$ws.Range("E2").Value = "  +3.91%  "
Set-TextValue "D3" "2.255.69"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "302.65"
$ws.Range("E5").Value = "  +3.36%  "
Set-TextValue "D6" "91.89"
$ws.Range("E6").Value = "  +6.06%  "
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.92%  "
Set-TextValue "D10" "54.07"
$ws.Range("E10").Value = "  +9.00%  "
Set-TextValue "D11" "32.04"
$ws.Range("E11").Value = "  +7.60%  "
Set-TextValue "D12" "0.0794"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("E14").Value = "  +3.62%  "
Set-TextValue "D15" "2.601.86"
$ws.Range("E15").Value = "  +2.56%  "
Set-TextValue "D16" "14.18"
$ws.Range("E16").Value = "  +4.34%  "
Set-TextValue "D17" "2.268.39"
$ws.Range("E17").Value = "  +0.43%  "
Set-TextValue "D18" "0.751"
$ws.Range("E18").Value = "  +4.48%  "
Set-TextValue "D19" "41.258.93"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("E20").Value = "  +9.29%  "
Set-TextValue "D21" "0.0₃0906"
$ws.Range("E21").Value = "  +3.12%  "
Set-TextValue "D22" "5.91"
$ws.Range("E22").Value = "  +3.46%  "
Set-TextValue "D23" "66.95"
$ws.Range("E23").Value = "  +3.16%  "
Set-TextValue "D24" "239.83"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +3.72%  "
Set-TextValue "D28" "23.71"
$ws.Range("E28").Value = "  +6.13%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.19"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D30" "9.64"
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D31" "34.09"
$ws.Range("E31").Value = "  +9.85%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "157.55"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +6.88%  "
Set-TextValue "D35" "0.0737"
$ws.Range("E35").Value = "  +5.06%  "
Set-TextValue "D36" "3.03"
$ws.Range("E36").Value = "  +9.04%  "
Set-TextValue "D37" "2.37"
$ws.Range("E37").Value = "  +1.58%  "
Set-TextValue "D38" "16.60"
$ws.Range("E38").Value = "  +9.90%  "
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("E41").Value = "  +6.48%  "
Set-TextValue "D42" "3.99"
$ws.Range("E42").Value = "  +7.53%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D43" "2.065.17"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "20.18"
$ws.Range("E44").Value = "  +17.39%  "
Set-TextValue "D45" "0.0277"
$ws.Range("E45").Value = "  +4.21%  "
Set-TextValue "D46" "10.11"
$ws.Range("E46").Value = "  +5.73%  "
Set-TextValue "D47" "2.96"
$ws.Range("E47").Value = "  +13.12%  "
Set-TextValue "D48" "2.07"
$ws.Range("E48").Value = "  -1.67%  "
Set-TextValue "D49" "2.473.59"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +3.90%  "
